$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting pan_number..benefit_amount to the right
$ws.Range("A1").EntireColumn.Insert()

$ws.Range("A1").Value = "employee_id"

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("G7").Select() | Out-Null
